$wb = $excel.ActiveWorkbook
$n = $wb.Worksheets.Count

# Step 1: rename all sheets to unique temporary names to avoid collisions
for ($i = 1; $i -le $n; $i++) {
    $wb.Worksheets.Item($i).Name = "__tmp_sheet_${i}__"
}

# Step 2: rename sheets to their final target names
$newNames = @("summ2", "summ38", "summ5", "summ0", "summ4", "summ11", "summ10", "summ3", "summ1")
for ($i = 1; $i -le $n; $i++) {
    $wb.Worksheets.Item($i).Name = $newNames[$i-1]
}

# Step 3: update data values, row 19 label, and delete row 22 for each sheet
# --- Sheet 1 (summ6 -> summ2) ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2,2).Value = [double]"-1.208511429276684"
$ws.Cells.Item(2,3).Value = [double]"0.3648035582433746"
$ws.Cells.Item(3,2).Value = [double]"-0.1803934923435788"
$ws.Cells.Item(3,3).Value = [double]"0.5325351500960895"
$ws.Cells.Item(4,2).Value = [double]"-0.5784152501010675"
$ws.Cells.Item(4,3).Value = [double]"0.006500402861944399"
$ws.Cells.Item(5,2).Value = [double]"-0.6349424010968406"
$ws.Cells.Item(5,3).Value = [double]"0.003923214549570322"
$ws.Cells.Item(6,2).Value = [double]"-0.3682913621529903"
$ws.Cells.Item(6,3).Value = [double]"0.2948592930609906"
$ws.Cells.Item(7,2).Value = [double]"0.2510694744637843"
$ws.Cells.Item(7,3).Value = [double]"0.08086893200165537"
$ws.Cells.Item(8,2).Value = [double]"0.00050712119548373"
$ws.Cells.Item(8,3).Value = [double]"1.559332920272018e-25"
$ws.Cells.Item(9,2).Value = [double]"0.008497014865138923"
$ws.Cells.Item(9,3).Value = [double]"0.07570601474320307"
$ws.Cells.Item(10,2).Value = [double]"-0.01058851679628893"
$ws.Cells.Item(10,3).Value = [double]"0.9386513274375926"
$ws.Cells.Item(11,2).Value = [double]"0.5947055596898402"
$ws.Cells.Item(11,3).Value = [double]"0.001792109610642233"
$ws.Cells.Item(12,2).Value = [double]"0.2341736922501507"
$ws.Cells.Item(12,3).Value = [double]"0.341927950917743"
$ws.Cells.Item(13,2).Value = [double]"-5.662736420757604e-05"
$ws.Cells.Item(13,3).Value = [double]"0.02806337491909065"
$ws.Cells.Item(14,2).Value = [double]"-5.82349080939675e-09"
$ws.Cells.Item(14,3).Value = [double]"0.8494718674062836"
$ws.Cells.Item(15,2).Value = [double]"0.01351156109207809"
$ws.Cells.Item(15,3).Value = [double]"0.7832071447911569"
$ws.Cells.Item(16,2).Value = [double]"0.08478756486935517"
$ws.Cells.Item(16,3).Value = [double]"0.01245272494398361"
$ws.Cells.Item(17,2).Value = [double]"1.396187054643594"
$ws.Cells.Item(17,3).Value = [double]"0.390333718179771"
$ws.Cells.Item(18,2).Value = [double]"-0.005580299018734783"
$ws.Cells.Item(18,3).Value = [double]"0.4660432272727042"
$ws.Cells.Item(19,1).Value = "street_length"
$ws.Cells.Item(19,2).Value = [double]"0.00307193215491784"
$ws.Cells.Item(19,3).Value = [double]"0.7182371163574114"
$ws.Cells.Item(20,2).Value = [double]"-0.313084704986604"
$ws.Cells.Item(20,3).Value = [double]"0.6540356810114093"
$ws.Cells.Item(21,2).Value = [double]"0.02146807917810796"
$ws.Cells.Item(21,3).Value = [double]"0.9708256071208975"
$ws.Rows.Item(22).Delete()

# --- Sheet 2 (summ8 -> summ38) ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2,2).Value = [double]"-1.545549697762445"
$ws.Cells.Item(2,3).Value = [double]"0.246668210082744"
$ws.Cells.Item(3,2).Value = [double]"-0.2138951607676538"
$ws.Cells.Item(3,3).Value = [double]"0.458307676897206"
$ws.Cells.Item(4,2).Value = [double]"-0.4899787503970911"
$ws.Cells.Item(4,3).Value = [double]"0.02340200492811008"
$ws.Cells.Item(5,2).Value = [double]"-0.5596434844623394"
$ws.Cells.Item(5,3).Value = [double]"0.01328493095682154"
$ws.Cells.Item(6,2).Value = [double]"-0.4568874663095736"
$ws.Cells.Item(6,3).Value = [double]"0.2010823251619885"
$ws.Cells.Item(7,2).Value = [double]"0.2758844545223895"
$ws.Cells.Item(7,3).Value = [double]"0.06542982007118747"
$ws.Cells.Item(8,2).Value = [double]"0.0004963757426945524"
$ws.Cells.Item(8,3).Value = [double]"2.911270792265136e-25"
$ws.Cells.Item(9,2).Value = [double]"0.009361237264078814"
$ws.Cells.Item(9,3).Value = [double]"0.04977111649248817"
$ws.Cells.Item(10,2).Value = [double]"-0.07519490516331377"
$ws.Cells.Item(10,3).Value = [double]"0.5857238988728422"
$ws.Cells.Item(11,2).Value = [double]"0.6265654630782985"
$ws.Cells.Item(11,3).Value = [double]"0.0008763001229296218"
$ws.Cells.Item(12,2).Value = [double]"0.2180561611236962"
$ws.Cells.Item(12,3).Value = [double]"0.369828427847853"
$ws.Cells.Item(13,2).Value = [double]"-6.198518673765901e-05"
$ws.Cells.Item(13,3).Value = [double]"0.01779746035463716"
$ws.Cells.Item(14,2).Value = [double]"2.096910741974832e-08"
$ws.Cells.Item(14,3).Value = [double]"0.5051785493642305"
$ws.Cells.Item(15,2).Value = [double]"0.06162421607380421"
$ws.Cells.Item(15,3).Value = [double]"0.2152174231589405"
$ws.Cells.Item(16,2).Value = [double]"0.1269710540146341"
$ws.Cells.Item(16,3).Value = [double]"0.0002224250056311896"
$ws.Cells.Item(17,2).Value = [double]"2.617994825712556"
$ws.Cells.Item(17,3).Value = [double]"0.10892297048577"
$ws.Cells.Item(18,2).Value = [double]"-0.002609140679519987"
$ws.Cells.Item(18,3).Value = [double]"0.729328100278245"
$ws.Cells.Item(19,1).Value = "street_length"
$ws.Cells.Item(19,2).Value = [double]"-0.0001840741800928264"
$ws.Cells.Item(19,3).Value = [double]"0.982770286550029"
$ws.Cells.Item(20,2).Value = [double]"-0.275398558521183"
$ws.Cells.Item(20,3).Value = [double]"0.6918376723653203"
$ws.Cells.Item(21,2).Value = [double]"-0.3353880778343419"
$ws.Cells.Item(21,3).Value = [double]"0.5727528389973318"
$ws.Rows.Item(22).Delete()

# --- Sheet 3 (summ1 -> summ5) ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2,2).Value = [double]"-1.948160670452912"
$ws.Cells.Item(2,3).Value = [double]"0.142307254730264"
$ws.Cells.Item(3,2).Value = [double]"0.06638858125458957"
$ws.Cells.Item(3,3).Value = [double]"0.8218979537024532"
$ws.Cells.Item(4,2).Value = [double]"-0.6135854437995872"
$ws.Cells.Item(4,3).Value = [double]"0.003537294161030785"
$ws.Cells.Item(5,2).Value = [double]"-0.6132766172303616"
$ws.Cells.Item(5,3).Value = [double]"0.005594690863110714"
$ws.Cells.Item(6,2).Value = [double]"-0.3811461332361865"
$ws.Cells.Item(6,3).Value = [double]"0.2853380392591771"
$ws.Cells.Item(7,2).Value = [double]"0.2014177694100888"
$ws.Cells.Item(7,3).Value = [double]"0.1585554376069737"
$ws.Cells.Item(8,2).Value = [double]"0.0005255047155013583"
$ws.Cells.Item(8,3).Value = [double]"9.867682049778345e-27"
$ws.Cells.Item(9,2).Value = [double]"0.007599933816155457"
$ws.Cells.Item(9,3).Value = [double]"0.1102422970100859"
$ws.Cells.Item(10,2).Value = [double]"0.02158463155946875"
$ws.Cells.Item(10,3).Value = [double]"0.875263168474696"
$ws.Cells.Item(11,2).Value = [double]"0.689536934285629"
$ws.Cells.Item(11,3).Value = [double]"0.0002305874798377344"
$ws.Cells.Item(12,2).Value = [double]"0.3378110280079548"
$ws.Cells.Item(12,3).Value = [double]"0.161986968716797"
$ws.Cells.Item(13,2).Value = [double]"-8.18264177571811e-05"
$ws.Cells.Item(13,3).Value = [double]"0.002232968131263127"
$ws.Cells.Item(14,2).Value = [double]"1.699527399895892e-08"
$ws.Cells.Item(14,3).Value = [double]"0.5871254050395139"
$ws.Cells.Item(15,2).Value = [double]"0.007081755889873285"
$ws.Cells.Item(15,3).Value = [double]"0.8861274070953923"
$ws.Cells.Item(16,2).Value = [double]"0.1085100480677959"
$ws.Cells.Item(16,3).Value = [double]"0.001291765291771929"
$ws.Cells.Item(17,2).Value = [double]"2.292762427678142"
$ws.Cells.Item(17,3).Value = [double]"0.152623561708769"
$ws.Cells.Item(18,2).Value = [double]"-0.0057410412624169"
$ws.Cells.Item(18,3).Value = [double]"0.4445565699226524"
$ws.Cells.Item(19,1).Value = "street_length"
$ws.Cells.Item(19,2).Value = [double]"0.005992291356254985"
$ws.Cells.Item(19,3).Value = [double]"0.4784435750785206"
$ws.Cells.Item(20,2).Value = [double]"0.6428702242549411"
$ws.Cells.Item(20,3).Value = [double]"0.3590214386909032"
$ws.Cells.Item(21,2).Value = [double]"-0.592944829505865"
$ws.Cells.Item(21,3).Value = [double]"0.3082767294008685"
$ws.Rows.Item(22).Delete()

# --- Sheet 4 (summ7 -> summ0) ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2,2).Value = [double]"-2.305779887516413"
$ws.Cells.Item(2,3).Value = [double]"0.08185053748082713"
$ws.Cells.Item(3,2).Value = [double]"-0.1568821247352234"
$ws.Cells.Item(3,3).Value = [double]"0.5831479098673089"
$ws.Cells.Item(4,2).Value = [double]"-0.5140151054485963"
$ws.Cells.Item(4,3).Value = [double]"0.01472070487481775"
$ws.Cells.Item(5,2).Value = [double]"-0.6463318256617641"
$ws.Cells.Item(5,3).Value = [double]"0.003444533001672051"
$ws.Cells.Item(6,2).Value = [double]"-0.5094100434221226"
$ws.Cells.Item(6,3).Value = [double]"0.133410114030642"
$ws.Cells.Item(7,2).Value = [double]"0.2568675474253251"
$ws.Cells.Item(7,3).Value = [double]"0.07074774811394302"
$ws.Cells.Item(8,2).Value = [double]"0.0004816323883960095"
$ws.Cells.Item(8,3).Value = [double]"7.76774073019689e-24"
$ws.Cells.Item(9,2).Value = [double]"0.01153798771994358"
$ws.Cells.Item(9,3).Value = [double]"0.01541809562575169"
$ws.Cells.Item(10,2).Value = [double]"0.06861089528040809"
$ws.Cells.Item(10,3).Value = [double]"0.6164160746706175"
$ws.Cells.Item(11,2).Value = [double]"0.6661713092575826"
$ws.Cells.Item(11,3).Value = [double]"0.000400944680832676"
$ws.Cells.Item(12,2).Value = [double]"0.196775945908252"
$ws.Cells.Item(12,3).Value = [double]"0.4155206526855978"
$ws.Cells.Item(13,2).Value = [double]"-5.004032818406614e-05"
$ws.Cells.Item(13,3).Value = [double]"0.05321565777843774"
$ws.Cells.Item(14,2).Value = [double]"-1.592143537836967e-09"
$ws.Cells.Item(14,3).Value = [double]"0.959183178494771"
$ws.Cells.Item(15,2).Value = [double]"0.02835893344959076"
$ws.Cells.Item(15,3).Value = [double]"0.5659665513130298"
$ws.Cells.Item(16,2).Value = [double]"0.1185266083839867"
$ws.Cells.Item(16,3).Value = [double]"0.0005959956718728316"
$ws.Cells.Item(17,2).Value = [double]"2.683812252658051"
$ws.Cells.Item(17,3).Value = [double]"0.09841310146664188"
$ws.Cells.Item(18,2).Value = [double]"-0.005731259207237257"
$ws.Cells.Item(18,3).Value = [double]"0.4504281304087165"
$ws.Cells.Item(19,1).Value = "street_length"
$ws.Cells.Item(19,2).Value = [double]"0.007538489074646924"
$ws.Cells.Item(19,3).Value = [double]"0.3719741096432607"
$ws.Cells.Item(20,2).Value = [double]"-0.02322459999329838"
$ws.Cells.Item(20,3).Value = [double]"0.9733911443526821"
$ws.Cells.Item(21,2).Value = [double]"-0.3838163629567677"
$ws.Cells.Item(21,3).Value = [double]"0.5100126932135297"
$ws.Rows.Item(22).Delete()

# --- Sheet 5 (summ20 -> summ4) ---
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2,2).Value = [double]"-2.243934949145964"
$ws.Cells.Item(2,3).Value = [double]"0.09547530100612443"
$ws.Cells.Item(3,2).Value = [double]"-0.2156487093821143"
$ws.Cells.Item(3,3).Value = [double]"0.4769379510786333"
$ws.Cells.Item(4,2).Value = [double]"-0.5382126424906842"
$ws.Cells.Item(4,3).Value = [double]"0.01537467817520871"
$ws.Cells.Item(5,2).Value = [double]"-0.5801908440186146"
$ws.Cells.Item(5,3).Value = [double]"0.01257793389895022"
$ws.Cells.Item(6,2).Value = [double]"-0.5370145280891087"
$ws.Cells.Item(6,3).Value = [double]"0.1335935431041417"
$ws.Cells.Item(7,2).Value = [double]"0.3300087002107472"
$ws.Cells.Item(7,3).Value = [double]"0.03565498979342356"
$ws.Cells.Item(8,2).Value = [double]"0.0004934308423529847"
$ws.Cells.Item(8,3).Value = [double]"2.408487451149946e-23"
$ws.Cells.Item(9,2).Value = [double]"0.009541472465294897"
$ws.Cells.Item(9,3).Value = [double]"0.0485809720847983"
$ws.Cells.Item(10,2).Value = [double]"0.04779634542267"
$ws.Cells.Item(10,3).Value = [double]"0.7310991565947573"
$ws.Cells.Item(11,2).Value = [double]"0.7547719664303377"
$ws.Cells.Item(11,3).Value = [double]"8.087839798146295e-05"
$ws.Cells.Item(12,2).Value = [double]"0.2998559437982312"
$ws.Cells.Item(12,3).Value = [double]"0.2210074649705169"
$ws.Cells.Item(13,2).Value = [double]"-4.633470181880416e-05"
$ws.Cells.Item(13,3).Value = [double]"0.07471044862035721"
$ws.Cells.Item(14,2).Value = [double]"3.162058426451745e-09"
$ws.Cells.Item(14,3).Value = [double]"0.9193026076753735"
$ws.Cells.Item(15,2).Value = [double]"0.04782179485002114"
$ws.Cells.Item(15,3).Value = [double]"0.3417016770653944"
$ws.Cells.Item(16,2).Value = [double]"0.1197530407439521"
$ws.Cells.Item(16,3).Value = [double]"0.0005678219369809006"
$ws.Cells.Item(17,2).Value = [double]"3.089268282773199"
$ws.Cells.Item(17,3).Value = [double]"0.05894278874746273"
$ws.Cells.Item(18,2).Value = [double]"-0.0007681723330367729"
$ws.Cells.Item(18,3).Value = [double]"0.920891541362383"
$ws.Cells.Item(19,1).Value = "street_length"
$ws.Cells.Item(19,2).Value = [double]"0.005136624863030022"
$ws.Cells.Item(19,3).Value = [double]"0.546575588383305"
$ws.Cells.Item(20,2).Value = [double]"-0.6468329615979476"
$ws.Cells.Item(20,3).Value = [double]"0.3565365541535152"
$ws.Cells.Item(21,2).Value = [double]"-0.5285223166130971"
$ws.Cells.Item(21,3).Value = [double]"0.364207090201348"
$ws.Rows.Item(22).Delete()

# --- Sheet 6 (summ5 -> summ11) ---
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2,2).Value = [double]"-1.95336926752722"
$ws.Cells.Item(2,3).Value = [double]"0.1378965293308488"
$ws.Cells.Item(3,2).Value = [double]"-0.1798241800364487"
$ws.Cells.Item(3,3).Value = [double]"0.5280711422521434"
$ws.Cells.Item(4,2).Value = [double]"-0.5330692136559815"
$ws.Cells.Item(4,3).Value = [double]"0.01090462852362926"
$ws.Cells.Item(5,2).Value = [double]"-0.6271568768293203"
$ws.Cells.Item(5,3).Value = [double]"0.004399253838408134"
$ws.Cells.Item(6,2).Value = [double]"-0.4990843753564941"
$ws.Cells.Item(6,3).Value = [double]"0.156986916215464"
$ws.Cells.Item(7,2).Value = [double]"0.2460773661611823"
$ws.Cells.Item(7,3).Value = [double]"0.08203181566374208"
$ws.Cells.Item(8,2).Value = [double]"0.0004935089674234283"
$ws.Cells.Item(8,3).Value = [double]"9.32094627301153e-25"
$ws.Cells.Item(9,2).Value = [double]"0.008912798198268025"
$ws.Cells.Item(9,3).Value = [double]"0.06067312280715074"
$ws.Cells.Item(10,2).Value = [double]"-0.03338715990997555"
$ws.Cells.Item(10,3).Value = [double]"0.8086590274763485"
$ws.Cells.Item(11,2).Value = [double]"0.6322509046612702"
$ws.Cells.Item(11,3).Value = [double]"0.0009384791583124212"
$ws.Cells.Item(12,2).Value = [double]"0.1986741953430593"
$ws.Cells.Item(12,3).Value = [double]"0.4138735964246042"
$ws.Cells.Item(13,2).Value = [double]"-4.683452880530393e-05"
$ws.Cells.Item(13,3).Value = [double]"0.07267073682577829"
$ws.Cells.Item(14,2).Value = [double]"1.365014754745949e-08"
$ws.Cells.Item(14,3).Value = [double]"0.6589453467491426"
$ws.Cells.Item(15,2).Value = [double]"0.01559151327908001"
$ws.Cells.Item(15,3).Value = [double]"0.7513853497326806"
$ws.Cells.Item(16,2).Value = [double]"0.1265438382361015"
$ws.Cells.Item(16,3).Value = [double]"0.0001969360929468159"
$ws.Cells.Item(17,2).Value = [double]"2.822685865663412"
$ws.Cells.Item(17,3).Value = [double]"0.08530296029321321"
$ws.Cells.Item(18,2).Value = [double]"-0.006595598471774791"
$ws.Cells.Item(18,3).Value = [double]"0.374574589474521"
$ws.Cells.Item(19,1).Value = "street_length"
$ws.Cells.Item(19,2).Value = [double]"0.005637536049226682"
$ws.Cells.Item(19,3).Value = [double]"0.5035637061853201"
$ws.Cells.Item(20,2).Value = [double]"-0.1085052229831149"
$ws.Cells.Item(20,3).Value = [double]"0.8746201556483444"
$ws.Cells.Item(21,2).Value = [double]"-0.3195730990422483"
$ws.Cells.Item(21,3).Value = [double]"0.5776856090415432"
$ws.Rows.Item(22).Delete()

# --- Sheet 7 (summ0 -> summ10) ---
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(2,2).Value = [double]"-1.264997950756201"
$ws.Cells.Item(2,3).Value = [double]"0.3457485026060095"
$ws.Cells.Item(3,2).Value = [double]"-0.2277536688022952"
$ws.Cells.Item(3,3).Value = [double]"0.4352136827678421"
$ws.Cells.Item(4,2).Value = [double]"-0.5653389847309567"
$ws.Cells.Item(4,3).Value = [double]"0.008563430280604665"
$ws.Cells.Item(5,2).Value = [double]"-0.5664021620236557"
$ws.Cells.Item(5,3).Value = [double]"0.01145367039812788"
$ws.Cells.Item(6,2).Value = [double]"-0.2874931259058655"
$ws.Cells.Item(6,3).Value = [double]"0.4137240576725468"
$ws.Cells.Item(7,2).Value = [double]"0.2904107558519265"
$ws.Cells.Item(7,3).Value = [double]"0.04952476652515633"
$ws.Cells.Item(8,2).Value = [double]"0.0005023900016445131"
$ws.Cells.Item(8,3).Value = [double]"3.949499849640381e-25"
$ws.Cells.Item(9,2).Value = [double]"0.009816219482805588"
$ws.Cells.Item(9,3).Value = [double]"0.04061744730341262"
$ws.Cells.Item(10,2).Value = [double]"-0.01716380448690352"
$ws.Cells.Item(10,3).Value = [double]"0.9007742899555718"
$ws.Cells.Item(11,2).Value = [double]"0.7639921229189819"
$ws.Cells.Item(11,3).Value = [double]"4.839766336128112e-05"
$ws.Cells.Item(12,2).Value = [double]"0.2781775905962146"
$ws.Cells.Item(12,3).Value = [double]"0.2462595986927041"
$ws.Cells.Item(13,2).Value = [double]"-6.042107880972548e-05"
$ws.Cells.Item(13,3).Value = [double]"0.0206561414070032"
$ws.Cells.Item(14,2).Value = [double]"2.256734252769587e-08"
$ws.Cells.Item(14,3).Value = [double]"0.4725831632448262"
$ws.Cells.Item(15,2).Value = [double]"0.05215576831785119"
$ws.Cells.Item(15,3).Value = [double]"0.2966112374856117"
$ws.Cells.Item(16,2).Value = [double]"0.1078679651295714"
$ws.Cells.Item(16,3).Value = [double]"0.001653711845330426"
$ws.Cells.Item(17,2).Value = [double]"1.959590797007318"
$ws.Cells.Item(17,3).Value = [double]"0.2248992184408927"
$ws.Cells.Item(18,2).Value = [double]"-0.008429788412275523"
$ws.Cells.Item(18,3).Value = [double]"0.2646596336983089"
$ws.Cells.Item(19,1).Value = "street_length"
$ws.Cells.Item(19,2).Value = [double]"-0.0006976587423623577"
$ws.Cells.Item(19,3).Value = [double]"0.9346037099590631"
$ws.Cells.Item(20,2).Value = [double]"-0.1079726628863033"
$ws.Cells.Item(20,3).Value = [double]"0.8753988134457655"
$ws.Cells.Item(21,2).Value = [double]"-0.727264537272208"
$ws.Cells.Item(21,3).Value = [double]"0.21090341962295"
$ws.Rows.Item(22).Delete()

# --- Sheet 8 (summ26 -> summ3) ---
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(2,2).Value = [double]"-1.490167937195062"
$ws.Cells.Item(2,3).Value = [double]"0.2670529938855229"
$ws.Cells.Item(3,2).Value = [double]"-0.0760620871508344"
$ws.Cells.Item(3,3).Value = [double]"0.796897760234079"
$ws.Cells.Item(4,2).Value = [double]"-0.5606642086537189"
$ws.Cells.Item(4,3).Value = [double]"0.008627184945741093"
$ws.Cells.Item(5,2).Value = [double]"-0.6001559648275961"
$ws.Cells.Item(5,3).Value = [double]"0.00786155330537655"
$ws.Cells.Item(6,2).Value = [double]"-0.2687336121869449"
$ws.Cells.Item(6,3).Value = [double]"0.4596384940715827"
$ws.Cells.Item(7,2).Value = [double]"0.2378422879254045"
$ws.Cells.Item(7,3).Value = [double]"0.1021319761878818"
$ws.Cells.Item(8,2).Value = [double]"0.0005096436420501998"
$ws.Cells.Item(8,3).Value = [double]"1.01092619396026e-24"
$ws.Cells.Item(9,2).Value = [double]"0.006785947129072043"
$ws.Cells.Item(9,3).Value = [double]"0.1560749981186389"
$ws.Cells.Item(10,2).Value = [double]"0.0662662671790585"
$ws.Cells.Item(10,3).Value = [double]"0.6324188589466588"
$ws.Cells.Item(11,2).Value = [double]"0.8199736888723063"
$ws.Cells.Item(11,3).Value = [double]"1.364153188673925e-05"
$ws.Cells.Item(12,2).Value = [double]"0.4531493824114969"
$ws.Cells.Item(12,3).Value = [double]"0.06137231096600126"
$ws.Cells.Item(13,2).Value = [double]"-4.828456185282547e-05"
$ws.Cells.Item(13,3).Value = [double]"0.06566074931631455"
$ws.Cells.Item(14,2).Value = [double]"-5.939126918558695e-09"
$ws.Cells.Item(14,3).Value = [double]"0.8509706091823134"
$ws.Cells.Item(15,2).Value = [double]"0.03741721573755293"
$ws.Cells.Item(15,3).Value = [double]"0.4545265725743622"
$ws.Cells.Item(16,2).Value = [double]"0.1013505453879778"
$ws.Cells.Item(16,3).Value = [double]"0.002904079602582348"
$ws.Cells.Item(17,2).Value = [double]"2.19995573043848"
$ws.Cells.Item(17,3).Value = [double]"0.1781226154013646"
$ws.Cells.Item(18,2).Value = [double]"-0.005635701623419848"
$ws.Cells.Item(18,3).Value = [double]"0.4621821649432107"
$ws.Cells.Item(19,1).Value = "street_length"
$ws.Cells.Item(19,2).Value = [double]"0.00256936366825809"
$ws.Cells.Item(19,3).Value = [double]"0.7648863320903267"
$ws.Cells.Item(20,2).Value = [double]"-0.2919133660542546"
$ws.Cells.Item(20,3).Value = [double]"0.6767144247766976"
$ws.Cells.Item(21,2).Value = [double]"-0.4553671243696893"
$ws.Cells.Item(21,3).Value = [double]"0.4445200831593223"
$ws.Rows.Item(22).Delete()

# --- Sheet 9 (summ2 -> summ1) ---
$ws = $wb.Worksheets.Item(9)
$ws.Cells.Item(2,2).Value = [double]"-1.774771277869069"
$ws.Cells.Item(2,3).Value = [double]"0.1840143428178622"
$ws.Cells.Item(3,2).Value = [double]"-0.07628563161282585"
$ws.Cells.Item(3,3).Value = [double]"0.8002407844520764"
$ws.Cells.Item(4,2).Value = [double]"-0.5032224055879881"
$ws.Cells.Item(4,3).Value = [double]"0.01962326657192942"
$ws.Cells.Item(5,2).Value = [double]"-0.5682126824695837"
$ws.Cells.Item(5,3).Value = [double]"0.01142660021156374"
$ws.Cells.Item(6,2).Value = [double]"-0.6141610795484449"
$ws.Cells.Item(6,3).Value = [double]"0.08259639660756"
$ws.Cells.Item(7,2).Value = [double]"0.2392260285884512"
$ws.Cells.Item(7,3).Value = [double]"0.1091017694149268"
$ws.Cells.Item(8,2).Value = [double]"0.0004975877138891932"
$ws.Cells.Item(8,3).Value = [double]"1.061940642976315e-24"
$ws.Cells.Item(9,2).Value = [double]"0.006383076434718389"
$ws.Cells.Item(9,3).Value = [double]"0.1806145667653013"
$ws.Cells.Item(10,2).Value = [double]"0.04547482014976555"
$ws.Cells.Item(10,3).Value = [double]"0.7414063254686658"
$ws.Cells.Item(11,2).Value = [double]"0.7044118756694404"
$ws.Cells.Item(11,3).Value = [double]"0.0002557849377449677"
$ws.Cells.Item(12,2).Value = [double]"0.3346243224105547"
$ws.Cells.Item(12,3).Value = [double]"0.1693684629783523"
$ws.Cells.Item(13,2).Value = [double]"-6.193253480182846e-05"
$ws.Cells.Item(13,3).Value = [double]"0.01728412050190954"
$ws.Cells.Item(14,2).Value = [double]"1.925155948545389e-08"
$ws.Cells.Item(14,3).Value = [double]"0.5381279763685627"
$ws.Cells.Item(15,2).Value = [double]"0.0415050881892997"
$ws.Cells.Item(15,3).Value = [double]"0.3980977885791498"
$ws.Cells.Item(16,2).Value = [double]"0.1125483932873322"
$ws.Cells.Item(16,3).Value = [double]"0.001323367343904416"
$ws.Cells.Item(17,2).Value = [double]"2.212043411321978"
$ws.Cells.Item(17,3).Value = [double]"0.1768638962496313"
$ws.Cells.Item(18,2).Value = [double]"-0.007870884757223917"
$ws.Cells.Item(18,3).Value = [double]"0.3072476913314025"
$ws.Cells.Item(19,1).Value = "street_length"
$ws.Cells.Item(19,2).Value = [double]"0.006233927052163545"
$ws.Cells.Item(19,3).Value = [double]"0.469497586919292"
$ws.Cells.Item(20,2).Value = [double]"-0.02002373735837542"
$ws.Cells.Item(20,3).Value = [double]"0.9772607240269219"
$ws.Cells.Item(21,2).Value = [double]"-0.6427685961227841"
$ws.Cells.Item(21,3).Value = [double]"0.2762971621642188"
$ws.Rows.Item(22).Delete()

